$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit moves the "consumption|district heat|District heat use" row (originally row 18)
# down to just after the "consumption|individual heat|individual heat use" row, i.e. it
# becomes the new row 22, while the rows that used to be 19-22 shift up to become 18-21.
# This is equivalent to: Cut row 18, Insert Cut Cells before (old) row 23.

# Capture the data that lives in row 18 before it is removed.
$vA = $ws.Range("A18").Value2
$vB = $ws.Range("B18").Value2
$vC = $ws.Range("C18").Value2
$vD = $ws.Range("D18").Value2
$vE = $ws.Range("E18").Value2
$vF = $ws.Range("F18").Value2
$vG = $ws.Range("G18").Value2
$vH = $ws.Range("H18").Value2

# Capture (and remove) the comment that sits on row 18 (E18) so it can be re-created
# on the row that now holds that data (new row 21, since the row that used to be E22
# shifts up to E21 once row 18 disappears).
$commentText = $null
$cmt = $ws.Range("E18").Comment
if ($cmt -ne $null) {
    $commentText = $cmt.Text()
    $cmt.Delete() | Out-Null
}

# Remove row 18 entirely; rows 19:50 shift up to become rows 18:49.
$ws.Rows.Item(18).Delete() | Out-Null

# Re-create the comment (that used to belong to the "individual heat use" row, old E22)
# on its new location E21.
if ($commentText -ne $null) {
    $ws.Range("E21").AddComment($commentText) | Out-Null
}

# Write the data that used to be in row 18 into the now-empty row 22.
$ws.Range("A22").Value2 = $vA
$ws.Range("B22").Value2 = $vB
$ws.Range("C22").Value2 = $vC
$ws.Range("D22").Value2 = $vD
$ws.Range("E22").Value2 = $vE
$ws.Range("F22").Value2 = $vF
$ws.Range("G22").Value2 = $vG
$ws.Range("H22").Value2 = $vH

# Move the selection to reflect where the edit finished (A18), matching the saved view.
$ws.Range("A18").Select() | Out-Null
